$d = $word.ActiveDocument

# Locate the paragraph that holds "Features " (the heading just added before
# the new paragraph we need to insert).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Features ") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Features ' paragraph"
}

# Collapse a range to the end of that paragraph (before its paragraph mark)
# and insert a brand new paragraph after it.
$endRange = $target.Range
$endRange.Collapse(0)  # wdCollapseEnd
$newParaRange = $endRange.InsertParagraphAfter()

# Re-fetch the paragraph that was just created so we can set its text and
# formatting precisely.
$newPara = $target.Next()
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.MoveEnd(1, -1) | Out-Null

$newPara.Range.Text = "Following are the features"

# Match formatting used for the sibling heading-style paragraph: bold,
# size 24 half-points (12pt), both-justified alignment.
$newPara.Range.Font.Bold = $true
$newPara.Range.Font.Size = 12
$newPara.Format.Alignment = 3  # wdAlignParagraphJustify
